# Weekly refresh of the "Fruta, Mercado Mayorista Lo Valledor de Santiago - Higo"
# consolidated sheet: the per-lot fields (Fecha, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Origen, Precio $/Kg) are reshuffled
# across the existing data rows (2-32). Every other column (ids, mercado, region,
# producto, categoria, variedad, unidad de comercializacion, kg/unidad) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as one "record" per row.
$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

$firstRow = 2
$lastRow = 32

# Snapshot every row's record BEFORE writing anything, so that later writes
# never read an already-overwritten source row.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rec = @{}
    foreach ($c in $cols) {
        $rec[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rec
}

# Target row -> source row (the record that ends up living at the target row).
$map = @{
    2 = 16; 3 = 30; 4 = 28; 5 = 25; 6 = 26; 7 = 3; 8 = 4; 9 = 27; 10 = 21;
    11 = 2; 12 = 10; 13 = 17; 14 = 18; 15 = 23; 16 = 13; 17 = 31; 18 = 32;
    19 = 5; 20 = 6; 21 = 7; 22 = 19; 23 = 20; 24 = 24; 25 = 29; 26 = 22;
    27 = 8; 28 = 9; 29 = 14; 30 = 15; 31 = 11; 32 = 12
}

foreach ($targetRow in ($map.Keys | Sort-Object)) {
    $sourceRow = $map[$targetRow]
    $rec = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value2 = $rec[$c]
    }
}
